$d = $word.ActiveDocument

# --- Locate and temporarily remove the hidden "_GoBack" bookmark ---
# The bookmark sits at the very end of the document, right after the
# final run of text in the last paragraph. We need to insert several
# new paragraphs there, before the bookmark. Re-adding a bookmark at
# exactly the same offset it started at trips an engine quirk, so we
# delete it now and re-create it later once the document has grown.
$d.Bookmarks.ShowHidden = $true
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# wdCollapseEnd = 0
$wdCollapseEnd = 0
# wdCharacter = 1
$wdCharacter = 1
# wdDarkYellow highlight index
$wdDarkYellow = 14
# wdNoHighlight
$wdNoHighlight = 0

# --- Paragraph: "Eine computerstimme..." already exists as the last paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# --- New paragraph 1: "clear_alias" (highlighted keyword) ---
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.Text = "clear_alias"
$p1Text = $p1.Range.Duplicate
$p1Text.MoveEnd($wdCharacter, -1)
$p1Text.HighlightColorIndex = $wdDarkYellow
$p1.Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.HighlightColorIndex = $wdNoHighlight

# --- New paragraph 2: description ---
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.Text = "löscht den Inhalt des Alias wieder"
$p2.Range.InsertParagraphAfter()

# --- New paragraph 3: "only_alias_num" (highlighted keyword) ---
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3.Range.Text = "only_alias_num"
$p3Text = $p3.Range.Duplicate
$p3Text.MoveEnd($wdCharacter, -1)
$p3Text.HighlightColorIndex = $wdDarkYellow
$p3.Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.HighlightColorIndex = $wdNoHighlight

# --- New paragraph 4: description, split into three runs around "ausgeführt" ---
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p4.Range.Text = "Code wird nur weiter ausgeführt wenn der im Alias gespeicherte Werte numerisch ist"
$p4.Range.InsertParagraphAfter()

# --- Final paragraph: re-home the bookmark here ---
$bmPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmPos = $bmPara.Range.Start
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
